$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet ("UMT-L2Manager Scenarios cred") after the existing "Roles" sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "UMT-L2Manager Scenarios cred"

# New test-case data.
$ws2.Range("A1").Value = "UserName"
$ws2.Range("A2").Value = "L2TEST11"
$ws2.Range("A3").Value = "L2TEST12"
$ws2.Range("A4").Value = "L2TEST13"

# Give the header cell the same bold / green-fill look used for the
# "ROLES" header on the Roles sheet.
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column width to match the source workbook's new sheet.
$ws2.Columns.Item(1).ColumnWidth = 16.67

# Match page setup used elsewhere in the workbook.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selection on the new sheet ends on A4.
$null = $ws2.Range("A4").Select()

# The new sheet becomes the active / selected tab.
$ws2.Activate()

Write-Output "done"
